# Auto-update of EPEX Spot prices, Gaz (gas) and CO2 workbook with the
# newest day of data (11-sep for spot prices, 2025-09-09 for gas/CO2).

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": append a new date column (CL) ------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting (bold header, border) of the previous header cell
# onto the new header cell before writing its value.
$wsSpot.Range("CK1").Copy()
$wsSpot.Range("CL1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSpot.Range("CL1").Value = "11-sep"
$wsSpot.Range("CL2").Value = 14
$wsSpot.Range("CL3").Value = 11.7
$wsSpot.Range("CL4").Value = 10.53
$wsSpot.Range("CL5").Value = 0
$wsSpot.Range("CL6").Value = 0
$wsSpot.Range("CL7").Value = 8.84
$wsSpot.Range("CL8").Value = 9.640000000000001
$wsSpot.Range("CL9").Value = 28.59
$wsSpot.Range("CL10").Value = 37.08
$wsSpot.Range("CL11").Value = 17.13
$wsSpot.Range("CL12").Value = 0.03
$wsSpot.Range("CL13").Value = -0.01
$wsSpot.Range("CL14").Value = -0.01
$wsSpot.Range("CL15").Value = -0.01
$wsSpot.Range("CL16").Value = -0.01
$wsSpot.Range("CL17").Value = -0.01
$wsSpot.Range("CL18").Value = -0.01
$wsSpot.Range("CL19").Value = 0
$wsSpot.Range("CL20").Value = 7.05
$wsSpot.Range("CL21").Value = 32.46
$wsSpot.Range("CL22").Value = 25.85
$wsSpot.Range("CL23").Value = 11.81
$wsSpot.Range("CL24").Value = 15.97
$wsSpot.Range("CL25").Value = 11.87

# --- Sheet "Gaz": append the newest daily price row ---------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A87").NumberFormat = "@"
$wsGaz.Range("A87").Value = "2025-09-09"
$wsGaz.Range("A87").Style = "Normal"
$wsGaz.Range("B87").Value = 32

# --- Sheet "CO2": append the newest daily price row ----------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A87").NumberFormat = "@"
$wsCO2.Range("A87").Value = "2025-09-09"
$wsCO2.Range("A87").Style = "Normal"
$wsCO2.Range("B87").Value = 75.8

Write-Host "Workbook updated: Prix Spot CL1:CL25, Gaz A87:B87, CO2 A87:B87"
